$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bracket")

$ws.Range("E4").Value = 'Wandering Albatross'
$ws.Range("M4").Value = 'Kob'
$ws.Range("D6").Value = 'Wandering Albatross'
$ws.Range("F8").Value = 'Stag'
$ws.Range("L8").Value = 'Kob'
$ws.Range("D10").Value = 'Stag'
$ws.Range("E12").Value = 'Stag'
$ws.Range("N14").Value = 'Cobra Lily'
$ws.Range("G16").Value = 'Giant Squid '
$ws.Range("K16").Value = 'Kob'
$ws.Range("N18").Value = 'Batfly'
$ws.Range("E20").Value = 'Tiger'
$ws.Range("N22").Value = 'Porcupine'
$ws.Range("F24").Value = 'Giant Squid '
$ws.Range("L24").Value = 'Fork-marked Lemur'
$ws.Range("E28").Value = 'Giant Squid '
$ws.Range("M28").Value = 'Northern Short-tailed Shrew'
$ws.Range("N30").Value = 'Parasitic Guest Ant'
$ws.Range("H32").Value = 'Giant Squid '
$ws.Range("I32").Value = 'Giant Squid '
$ws.Range("J32").Value = 'Northern Elephant Seal'
$ws.Range("C35").Value = 'Sparklemuffin Peacock Spider'
$ws.Range("F40").Value = 'Painted Wild Dog'
$ws.Range("N42").Value = 'Pronghorn'
$ws.Range("E44").Value = 'Madagascan Sunset Moth'
$ws.Range("M44").Value = 'Leatherback Turtle'
$ws.Range("G48").Value = 'Painted Wild Dog'
$ws.Range("E52").Value = 'Flat Lizard'
$ws.Range("M52").Value = 'Great White Shark'
$ws.Range("D54").Value = 'Flat Lizard'
$ws.Range("N54").Value = 'Great White Shark'
$ws.Range("L56").Value = 'Forest Elephant'
$ws.Range("D58").Value = 'Golden-headed Lion Tamarin '
$ws.Range("M60").Value = 'Forest Elephant'
$ws.Range("N62").Value = 'Forest Elephant'
